$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C header (date label), same style as B1 (bold/border/center)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "13-01-2023"

# Reorder rows: funds move to rows 2-10, avg/total move to rows 11-12
# Row data after edit: (A, B, C)
$data = @(
    @("1822 Raices Valores Negociables", 27202.22, 26329.59),
    @("Alpha Acciones", 370309.86, 368466.44),
    @("Alpha Mega", 538018.11, 534458.55),
    @("Alpha Recursos Naturales", 268255.2, 269384.72),
    @("Alpha renta balan global", 124900.04, 118364.76),
    @("Fima Acciones", 256290.21, 253729.65),
    @("Fima PB Acciones", 301155.97, 303089.32),
    @("HF Acciones Argentinas", 9540.9, 0),
    @("HF Acciones Lideres", 138994.33, 139293),
    @("avg", 226074.09, 223679.56),
    @("total", 2034666.84, 2013116.03)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
